$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 30: add E30
$ws.Range("E30").Value = 18

# Row 31: add C31, D31, E31
$ws.Range("C31").Value = 10
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 18

# Row 32: add C32, D32, E32
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0

# Row 33: add C33, D33, E33
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0

# Row 34: add C34, D34, E34
$ws.Range("C34").Value = 10
$ws.Range("D34").Value = 0.25
$ws.Range("E34").Value = 17

# Row 35: add C35, D35, E35
$ws.Range("C35").Value = 15
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 18

# Row 36: add C36, D36 (no E36)
$ws.Range("C36").Value = 8
$ws.Range("D36").Value = 0

# New cell H2 with formula
$ws.Range("H2").Formula = "=SUM(F2:F123)"

# Update selection to E36
$ws.Range("E36").Select()
